$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price / volume updates ---
$ws.Range("D2").Value = "42.652.96"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "2.267.87"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'249.45"

$ws.Range("D6").Value = "'0.641"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("D7").Value = "'76.68"
$ws.Range("E7").Value = "  +6.58%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -3.86%  "

$ws.Range("D10").Value = "'39.76"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("D11").Value = "'0.0967"

$ws.Range("D12").Value = "'7.27"
$ws.Range("E12").Value = "  -2.36%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("D14").Value = "2.605.41"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "'14.95"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("D16").Value = "'0.860"
$ws.Range("E16").Value = "  -2.93%  "

$ws.Range("D17").Value = "2.263.65"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").Value = "42.562.39"
$ws.Range("E18").Value = "  -0.46%  "

$ws.Range("D19").Value = "0.0₃0987"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").Value = "'6.16"

$ws.Range("D21").Value = "'72.02"
$ws.Range("E21").Value = "  -1.46%  "

$ws.Range("D22").Value = "'235.12"
$ws.Range("E22").Value = "  +0.72%  "

$ws.Range("D23").Value = "'2.13"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D26").Value = "'11.24"
$ws.Range("E26").Value = "  -1.13%  "

$ws.Range("D27").Value = "'2.36"
$ws.Range("E27").Value = "  -3.08%  "

$ws.Range("D28").Value = "'2.17"
$ws.Range("E28").Value = "  -1.42%  "

$ws.Range("D29").Value = "'167.30"
$ws.Range("E29").Value = "  -0.15%  "

$ws.Range("D30").Value = "'20.84"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("D31").Value = "'6.39"
$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("D32").Value = "'0.0852"
$ws.Range("E32").Value = "  +5.72%  "

$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "  -4.02%  "

$ws.Range("D34").Value = "'30.72"
$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").Value = "'4.55"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("D38").Value = "'0.0303"
$ws.Range("E38").Value = "  -3.20%  "

$ws.Range("D39").Value = "'13.72"
$ws.Range("E39").Value = "  +7.09%  "

$ws.Range("D40").Value = "'2.25"
$ws.Range("E40").Value = "  -3.33%  "

$ws.Range("D41").Value = "'5.83"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").Value = "'109.06"

$ws.Range("D44").Value = "'60.92"
$ws.Range("E44").Value = "  -1.97%  "

$ws.Range("D45").Value = "'8.83"
$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("E49").Value = "  -2.72%  "

$ws.Range("D50").Value = "'1.16"
$ws.Range("E50").Value = "  -2.97%  "

# --- Row 24/25: Dai <-> WEMIXToken reorder + updated values ---
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D24").Value = "'3.79"
$ws.Range("E24").Value = "  -5.63%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.11%  "

# --- Row 46/47: FTXToken <-> Cronos reorder + updated values ---
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.0999"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").Value = "'4.62"
$ws.Range("E47").Value = "  -8.61%  "

# --- Row 51: Bonk -> SynthetixNetwork ---
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'4.20"
$ws.Range("E51").Value = "  -1.56%  "

